$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 177
$ws.Range("I5").Value = 133.75
$ws.Range("K5").Value = 133.75
$ws.Range("M5").Value = -18.75
$ws.Range("H69").Value = 1215
$ws.Range("I69").Value = 1215
$ws.Range("K69").Value = 3645
$ws.Range("M69").Value = -2771
$ws.Range("H72").Value = 1215
$ws.Range("I72").Value = 1215
$ws.Range("K72").Value = 10935
$ws.Range("M72").Value = -6567
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H97").Value = 3392.1667
$ws.Range("J97").Value = 3392.1667
$ws.Range("L97").Value = 10176.5001
$ws.Range("N97").Value = -11168.5001
$ws.Range("H98").Value = 493.45456
$ws.Range("I98").Value = 472.375
$ws.Range("J98").Value = 549.6667
$ws.Range("K98").Value = 472.375
$ws.Range("L98").Value = 549.6667
$ws.Range("M98").Value = 1025.625
$ws.Range("N98").Value = -3545.6667
$ws.Range("H99").Value = 780.4
$ws.Range("I99").Value = 756.75
$ws.Range("J99").Value = 875
$ws.Range("K99").Value = 2270.25
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = -772.25
$ws.Range("N99").Value = -5621
$ws.Range("H122").Value = 493.45456
$ws.Range("I122").Value = 472.375
$ws.Range("J122").Value = 549.6667
$ws.Range("K122").Value = 1417.125
$ws.Range("L122").Value = 1649.0001
$ws.Range("M122").Value = 1032.875
$ws.Range("N122").Value = -6549.0001
$ws.Range("H137").Value = 2015.5714
$ws.Range("I137").Value = 1865.5
$ws.Range("J137").Value = 2215.6667
$ws.Range("K137").Value = 5596.5
$ws.Range("L137").Value = 6647.000100000001
$ws.Range("M137").Value = -3046.5
$ws.Range("N137").Value = -11747.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4421
$ws.Range("I2").Value = 4540.2856
$ws.Range("J2").Value = 4212.25
$ws.Range("K2").Value = 4540.2856
$ws.Range("L2").Value = 4212.25
$ws.Range("M2").Value = -4427.2856
$ws.Range("N2").Value = -4438.25
$ws.Range("H61").Value = 2277.8572
$ws.Range("I61").Value = 2028.1578
$ws.Range("K61").Value = 2028.1578
$ws.Range("M61").Value = -1816.1578
$ws.Range("H116").Value = 4421
$ws.Range("I116").Value = 4540.2856
$ws.Range("J116").Value = 4212.25
$ws.Range("K116").Value = 4540.2856
$ws.Range("L116").Value = 4212.25
$ws.Range("M116").Value = -2246.2856
$ws.Range("N116").Value = -8800.25
$ws.Range("H136").Value = 2277.8572
$ws.Range("I136").Value = 2028.1578
$ws.Range("K136").Value = 6084.4734
$ws.Range("M136").Value = -3534.4734

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4421
$ws.Range("I3").Value = 4540.2856
$ws.Range("J3").Value = 4212.25
$ws.Range("K3").Value = 4540.2856
$ws.Range("L3").Value = 4212.25
$ws.Range("M3").Value = -4426.2856
$ws.Range("N3").Value = -4440.25
$ws.Range("H20").Value = 4249.5
$ws.Range("I20").Value = 4249.5
$ws.Range("K20").Value = 4249.5
$ws.Range("M20").Value = -4002.5
$ws.Range("H134").Value = 4324.0454
$ws.Range("I134").Value = 1515.1875
$ws.Range("K134").Value = 4545.5625
$ws.Range("M134").Value = -2010.5625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 12999.333
$ws.Range("J23").Value = 12999.333
$ws.Range("L23").Value = 12999.333
$ws.Range("N23").Value = -13479.333
$ws.Range("H27").Value = 12999.333
$ws.Range("J27").Value = 12999.333
$ws.Range("L27").Value = 12999.333
$ws.Range("N27").Value = -13383.333
$ws.Range("H86").Value = 11799.8
$ws.Range("I86").Value = 9000
$ws.Range("J86").Value = 12499.75
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 12499.75
$ws.Range("M86").Value = -7877
$ws.Range("N86").Value = -14745.75
$ws.Range("H89").Value = 11799.8
$ws.Range("I89").Value = 9000
$ws.Range("J89").Value = 12499.75
$ws.Range("K89").Value = 45000
$ws.Range("L89").Value = 62498.75
$ws.Range("M89").Value = -39384
$ws.Range("N89").Value = -73730.75
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1239.8
$ws.Range("I68").Value = 1066.6666
$ws.Range("J68").Value = 1499.5
$ws.Range("K68").Value = 3199.9998
$ws.Range("L68").Value = 4498.5
$ws.Range("M68").Value = -2388.9998
$ws.Range("N68").Value = -6120.5
$ws.Range("H71").Value = 1239.8
$ws.Range("I71").Value = 1066.6666
$ws.Range("J71").Value = 1499.5
$ws.Range("K71").Value = 9599.999400000001
$ws.Range("L71").Value = 13495.5
$ws.Range("M71").Value = -5543.999400000001
$ws.Range("N71").Value = -21607.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H102").Value = 4780.3335
$ws.Range("I102").Value = 5371
$ws.Range("J102").Value = 3599
$ws.Range("K102").Value = 5371
$ws.Range("L102").Value = 3599
$ws.Range("M102").Value = -3749
$ws.Range("N102").Value = -6843
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 4725.5
$ws.Range("I132").Value = 3876.5
$ws.Range("K132").Value = 11629.5
$ws.Range("M132").Value = -9099.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3343003
$ws.Range("I5").Value = 9009
$ws.Range("K5").Value = 9009
$ws.Range("M5").Value = -8897
$ws.Range("H126").Value = 4070
$ws.Range("I126").Value = 3105
$ws.Range("K126").Value = 9315
$ws.Range("M126").Value = -6845
$ws.Range("H136").Value = 3232.3333
$ws.Range("I136").Value = 2679
$ws.Range("K136").Value = 8037
$ws.Range("M136").Value = -5487
